# Append the 2020-06-12 row (row 93) to the COVID data table on the
# active sheet, mirroring what Excel does when you type a new row of
# data directly below an existing Excel Table ("Tabela1"): the table
# auto-expands, the new row inherits formatting from the row above it,
# and the sheet's used-range/selection move down to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the table (Tabela1) by one row: A1:J92 -> A1:J93
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the previous last data row (89) onto the new
# row so the new row picks up the correct number formats/borders
# instead of being left completely unformatted.
$ws.Range("A89:J89").Copy()
$ws.Range("A93:J93").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Cells.Item(93, 1).Value = 43993
$ws.Cells.Item(93, 2).Value = 86328
$ws.Cells.Item(93, 3).Value = 702
$ws.Cells.Item(93, 4).Value = 1490
$ws.Cells.Item(93, 5).Value = 2
$ws.Cells.Item(93, 6).Value = 6
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 109
$ws.Cells.Item(93, 10).Value = 0

# Move the selection to the newly added row, matching the saved
# worksheet's cursor position.
$ws.Range("A93:J93").Select() | Out-Null
